$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Grab the formatting (empty run + bold run) of the existing
#    "Meta description" paragraph (2nd paragraph in the doc) so we can
#    reuse the same run-split for the new "Play Cosmic Jewels Online
#    Slot Game for Free" paragraph that needs to appear near the end
#    of the document, right before the closing "feature image" prompt.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$boldRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start)
$boldRange.Collapse(1)
$boldRange.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$labelRange = $d.Range($metaPara.Range.Start, $boldRange.End)
$savedFormattedText = $labelRange.FormattedText

# ------------------------------------------------------------------
# 2. Insert a brand-new paragraph right after "No progressive jackpot"
#    (i.e. immediately before the last paragraph) and stamp it with the
#    copied empty-run/bold-run formatting, then swap its text.
# ------------------------------------------------------------------
$noProgPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$noProgPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$newPara.Style = "Normal"
$newPara.Range.FormattedText = $savedFormattedText

$newPara2 = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$newPara2.Range.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, 1, $false, "Play Cosmic Jewels Online Slot Game for Free", 2) | Out-Null

# ------------------------------------------------------------------
# 3. Swap the text of the last paragraph (the italic image prompt) for
#    the meta-description copy, keeping its italic run intact.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.Find.Execute("Create an eye-catching feature image for Cosmic Jewels that captures the excitement and adventure of the game. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The Maya warrior should be holding a space helmet in one hand and a shining diamond in the other, with stars and planets in the background. Use bold colors and dynamic poses to make the image stand out and draw players in. Make sure the image accurately represents the space adventure theme of the game and entices players to join in the hunt for valuable jewels.", $true, $false, $false, $false, $false, $true, 1, $false, "Explore outer space and earn high payouts with Cosmic Jewels. Play for free, discover wilds, Scatter Bucks, free games and a good RTP.", 2) | Out-Null

# ------------------------------------------------------------------
# 4. Finally, remove the original "Meta description" paragraph from
#    the top of the document.
# ------------------------------------------------------------------
$metaPara2 = $d.Paragraphs.Item(2)
$metaPara2.Range.Delete()

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
